$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.565.71'
$ws.Range('E2').Value = '  -1.51%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.353.87'
$ws.Range('E3').Value = '  -1.25%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.35'
$ws.Range('E5').Value = '  -1.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '108.38'
$ws.Range('E6').Value = '  +2.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.636'
$ws.Range('E7').Value = '  -1.51%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.619'
$ws.Range('E9').Value = '  -5.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.72'
$ws.Range('E10').Value = '  +0.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0925'
$ws.Range('E11').Value = '  -1.76%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.49'
$ws.Range('E12').Value = '  -1.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.998'
$ws.Range('E13').Value = '  -2.52%  '
$ws.Range('E14').Value = '  +0.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.95'
$ws.Range('E15').Value = '  -7.56%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.710.18'
$ws.Range('E16').Value = '  -1.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.359.92'
$ws.Range('E17').Value = '  -1.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.686.99'
$ws.Range('E18').Value = '  -1.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.74'
$ws.Range('E19').Value = '  +0.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000106'
$ws.Range('E20').Value = '  -2.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '76.66'
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.64'
$ws.Range('E22').Value = '  +6.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '256.78'
$ws.Range('E23').Value = '  -7.20%  '
$ws.Range('E24').Value = '  -3.71%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.35'
$ws.Range('E25').Value = '  -3.37%  '
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.44'
$ws.Range('E27').Value = '  -2.40%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '22.88'
$ws.Range('E28').Value = '  -0.16%  '
$ws.Range('E29').Value = '  +3.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '174.84'
$ws.Range('E30').Value = '  -0.86%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '36.73'
$ws.Range('E31').Value = '  -2.96%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0890'
$ws.Range('E32').Value = '  -4.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.11'
$ws.Range('E33').Value = '  +3.40%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.91'
$ws.Range('E34').Value = '  -9.25%  '
$ws.Range('E35').Value = '  +18.94%  '
$ws.Range('E36').Value = '  -1.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.62'
$ws.Range('E37').Value = '  -4.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0362'
$ws.Range('E38').Value = '  -1.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.82'
$ws.Range('E39').Value = '  -7.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.66'
$ws.Range('E40').Value = '  -5.67%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.238'
$ws.Range('E41').Value = '  +2.61%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '72.01'
$ws.Range('E42').Value = '  +4.16%  '
$ws.Range('E43').Value = '  -7.14%  '
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('B45').Value = 'Celestia'
$ws.Range('C45').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '11.98'
$ws.Range('E45').Value = '  -3.99%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '113.29'
$ws.Range('E46').Value = '  -8.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.49'
$ws.Range('E47').Value = '  -2.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.13'
$ws.Range('E48').Value = '  -4.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '84.04'
$ws.Range('E49').Value = '  -11.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '74.07'
$ws.Range('E50').Value = '  +1.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.27'
$ws.Range('E51').Value = '  -3.52%  '
